# Rename the "_old"/"_new" suffixed header columns to use the respective
# format-version names ("_FV2210" for the old/source file, "_FV2304" for the
# new/target file), then turn the header row + data range into a proper
# Excel Table (ListObject) and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldHeaders = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)

$newHeaders = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

# Columns A-J (1-10): "_old" -> "_FV2210"
for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i]
}

# Column K (11) stays "diff"

# Columns L-U (12-21): "_new" -> "_FV2304"
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newHeaders[$i]
}

# Turn A1:U58 into a native Excel Table, reusing the header row text above as
# the column names, with an autofilter on the header row.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U58"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

# Freeze the header row (split after row 1, top-left cell of the scrolling
# pane is A2).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
